# "Updated: ut 30. 08. 2022" data refresh for the Slovakia COVID daily-stats
# sheet: a handful of already-recorded AgTests/AgPosit/hospitalised figures
# (columns F/G/H) were corrected for late-arriving data, and three new days
# (rows 905-907) plus a brand-new day (row 908, 2022-08-29) were filled in.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column F (AgTests) corrections ---------------------------------------
$ws.Range("F771").Value = 9520
$ws.Range("F776").Value = 15184
$ws.Range("F778").Value = 9312
$ws.Range("F782").Value = 11051
$ws.Range("F785").Value = 7341
$ws.Range("F789").Value = 8333
$ws.Range("F796").Value = 4805
$ws.Range("F803").Value = 3986
$ws.Range("F810").Value = 4117
$ws.Range("F831").Value = 2304
$ws.Range("F845").Value = 4112

# --- Column H (Pocet.hospitalizovanych) corrections, rows 861-874 ---------
$ws.Range("H861").Value = 379
$ws.Range("H862").Value = 397
$ws.Range("H863").Value = 388
$ws.Range("H864").Value = 426
$ws.Range("H865").Value = 476
$ws.Range("H866").Value = 502
$ws.Range("H867").Value = 535
$ws.Range("H868").Value = 547
$ws.Range("H869").Value = 566
$ws.Range("H870").Value = 592
$ws.Range("H871").Value = 643
$ws.Range("H872").Value = 699
$ws.Range("H873").Value = 722
$ws.Range("H874").Value = 743

# --- Column H corrections, rows 880-894 ------------------------------------
$ws.Range("H880").Value = 721
$ws.Range("H881").Value = 696
$ws.Range("H882").Value = 652
$ws.Range("H883").Value = 656
$ws.Range("H884").Value = 583
$ws.Range("H885").Value = 609
$ws.Range("H886").Value = 647
$ws.Range("H888").Value = 558
$ws.Range("H889").Value = 542
$ws.Range("H890").Value = 523
$ws.Range("H891").Value = 484
$ws.Range("H893").Value = 533
$ws.Range("H894").Value = 507

# --- Rows 895-900: mix of F/G/H corrections --------------------------------
$ws.Range("F895").Value = 3155
$ws.Range("H895").Value = 504

$ws.Range("H896").Value = 499

$ws.Range("F897").Value = 3531
$ws.Range("G897").Value = 285
$ws.Range("H897").Value = 505

$ws.Range("H898").Value = 475

$ws.Range("F899").Value = 881
$ws.Range("G899").Value = 79
$ws.Range("H899").Value = 491

$ws.Range("F900").Value = 1168
$ws.Range("G900").Value = 114
$ws.Range("H900").Value = 521

# --- Rows 901-904: F/G corrections plus newly-filled H column -------------
$ws.Range("F901").Value = 4781
$ws.Range("G901").Value = 449
$ws.Range("H901").Value = 516

$ws.Range("F902").Value = 2853
$ws.Range("G902").Value = 310
$ws.Range("H902").Value = 500

$ws.Range("F903").Value = 3047
$ws.Range("G903").Value = 300
$ws.Range("H903").Value = 473

$ws.Range("F904").Value = 5090
$ws.Range("G904").Value = 491
$ws.Range("H904").Value = 478

# --- Rows 905-907 previously only had A:E filled in; fill F/G/H now -------
$ws.Range("F905").Value = 4261
$ws.Range("G905").Value = 216
$ws.Range("H905").Value = 483

$ws.Range("F906").Value = 487
$ws.Range("G906").Value = 42
$ws.Range("H906").Value = 479

$ws.Range("F907").Value = 507
$ws.Range("G907").Value = 58
$ws.Range("H907").Value = 481

# --- New row 908 (2022-08-29) ----------------------------------------------
$ws.Range("A908").Value = 44802
$ws.Range("B908").Value = 1833873
$ws.Range("C908").Value = 488
$ws.Range("D908").Value = 125
$ws.Range("E908").Value = 20349
$ws.Range("F908").Value = 397
$ws.Range("G908").Value = 44
$ws.Range("H908").Value = 481
